# "Money S4 almost done"
# Adds an "AdapterId" column to the two small lookup tables (Tabulka18 /
# Tabulka29) on the "empty start" sheet, a "ChangedOn" column to the
# ExternalId/AdapterId/AccountId table (Tabulka5711), and fills in a data
# row for Tabulka510, mirroring the pattern already used on sheet "List1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("empty start")

# ---------------------------------------------------------------------
# Table "Tabulka29" (F2:I5 ID/Title/AccountId/ExternalId) needs to shift
# one column to the right (G2:K5) to make room, and gain an "AdapterId"
# column. Move the existing values/styles first, then resize the table.
# ---------------------------------------------------------------------
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4104) | Out-Null   # xlPasteAll
$ws.Range("F1").ClearContents()

$ws.Range("F2:I2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4104) | Out-Null
$ws.Range("F2:I2").ClearContents()

$ws.Range("H3").Copy() | Out-Null
$ws.Range("I3").PasteSpecial(-4104) | Out-Null
$ws.Range("H3").ClearContents()
$ws.Range("H3").ClearFormats()

$ws.Range("F4:H4").Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4104) | Out-Null
$ws.Range("F4:H4").ClearContents()
$ws.Range("F4").ClearFormats()

$lo29 = $ws.ListObjects.Item("Tabulka29")
$lo29.Resize($ws.Range("G2:K5"))
$ws.Range("K2").Value = "AdapterId"
$ws.Range("J3").Value = "AA1"
$ws.Range("K4").Value = 1

# ---------------------------------------------------------------------
# Table "Tabulka18" (A2:D5 ID/Title/AccountId/ExternalId) simply gains an
# "AdapterId" column in place (E2:E5).
# ---------------------------------------------------------------------
$lo18 = $ws.ListObjects.Item("Tabulka18")
$lo18.Resize($ws.Range("A2:E5"))
$ws.Range("E2").Value = "AdapterId"
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 1

# ---------------------------------------------------------------------
# Table "Tabulka510" (A9:D12 ID/Title/ExternalId/AdapterId) gets its
# first data row filled in.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Test 1"
$ws.Range("C10").Value = "A1"
$ws.Range("D10").Value = 1

# ---------------------------------------------------------------------
# Table "Tabulka5711" (A15:C19 ExternalId/AdapterId/AccountId) gains a
# "ChangedOn" column (D15:D19), and the previously-empty AdapterId data
# cells (B16:B18) now carry values.
# ---------------------------------------------------------------------
$lo5711 = $ws.ListObjects.Item("Tabulka5711")
$lo5711.Resize($ws.Range("A15:D19"))
$ws.Range("D15").Value = "ChangedOn"

$ws.Range("A16").Value = "A1"
$ws.Range("B16").Value = 1
$ws.Range("A17").Value = "A2"
$ws.Range("B17").Value = 1
$ws.Range("A18").Value = "A3"
$ws.Range("B18").Value = 1

$ws.Range("C16").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("C18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("C19").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Cosmetic: the last selection the author left behind.
# ---------------------------------------------------------------------
$ws.Range("G15").Select()
